$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A86").Value = "GRT-USD"
